# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (column G) values for rows 2-23, replacing the previous Strike# derived values.
$kValues = @(5, 3, 9, 6, 3, 7, 3, 2, 4, 8, 7, 6, 1, 7, 4, 1, 9, 4, 7, 10, 3, 2)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}

$wb.Save()
